$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 1.1
$ws.Range("C2").Value = 0.97
$ws.Range("D2").Value = 0.83
$ws.Range("E2").Value = 1.85
$ws.Range("G2").Value = 2.11

# Row 3 updates
$ws.Range("D3").Value = 0.83
$ws.Range("E3").Value = 1.85

# Row 4 updates
$ws.Range("D4").Value = 0.83
$ws.Range("E4").Value = 1.85
$ws.Range("F4").Value = 0.9399999999999999
$ws.Range("G4").Value = 0.96
$ws.Range("I4").Value = 0.7
